{"js": "// Add a new \"Author\" styled paragraph right after the \"Edison Achalma\"\n// author paragraph, containing the author's affiliation.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the paragraph styled \"Author\" whose text is exactly \"Edison Achalma\".\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Author\" && para.text.trim() === \"Edison Achalma\") {\n    target = para;\n    break;\n  }\n}\n\nif (target) {\n  // Insert a paragraph break right after the existing text, followed by the\n  // new affiliation text, in a single insertText call. Embedding the break\n  // inside the text that is appended to the end of the \"Edison Achalma\"\n  // paragraph makes the new paragraph automatically inherit the same\n  // \"Author\" style, and keeps \"Edison Achalma\" itself untouched.\n  const endRange = target.getRange(\"End\");\n  endRange.insertText(\n    \"\\rEscuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# Add a new \"Author\" styled paragraph right after the \"Edison Achalma\"\n# author paragraph, containing the author's affiliation.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph styled \"Author\" whose text is exactly \"Edison Achalma\".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Author\" -and $p.Range.Text.Trim() -eq \"Edison Achalma\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Insert a paragraph break right after the existing text, followed by the\n    # new affiliation text. Because the break is embedded inside the run that\n    # is appended to the end of the \"Edison Achalma\" paragraph, the new\n    # paragraph automatically inherits the same \"Author\" style.\n    $target.Range.InsertAfter(\"`rEscuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\")\n}\n"}
